# Applies the "documented-system" StructureDefinition metadata refresh:
#   - Version bump 5.0.0 -> 6.0.0
#   - Date refresh
#   - Publisher / Contact -> Publisher + Jurisdiction rework (drops the
#     duplicated "Contact" row entirely)
#   - Elements sheet: the root Extension's Short/Definition cells get the
#     resource-specific text instead of the generic Extension placeholder.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (row 9) - previously blank
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be "Contact" / "No display for ContactDetail";
# it becomes "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -
# remove it entirely (rows below shift up)
$meta.Rows("11").Delete()

$elements = $wb.Worksheets("Elements")

# Root Extension row (row 2): Short / Definition columns (K/L) get the
# resource-specific description instead of the generic placeholder text
$elements.Range("K2").Value = "Documented System"
$elements.Range("L2").Value = "Source system or jurisdiction of the legal document"
